# Insert a new data row at row 189 (pushing the existing rows 189..295 down
# to 190..296) and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 189 and below down by one row.
$ws.Rows("189:189").Insert()

# Populate the newly inserted row 189 with its data.
$ws.Range("A189").Value = 9
$ws.Range("B189").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C189").Value = "Metropolitana"
$ws.Range("D189").Value = 44582
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = 100112044
$ws.Range("G189").Value = "Perejil"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 79
$ws.Range("K189").Value = 13000
$ws.Range("L189").Value = 14000
$ws.Range("M189").Value = 13494
$ws.Range("N189").Value = "`$/docena de atados"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 4498
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"
